$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = "sd"
$ws.Range("J5").Value = "Statement-non-opinion"
$ws.Range("I29").Value = "sd"
$ws.Range("J29").Value = "Statement-non-opinion"
$ws.Range("I32").Value = "sd"
$ws.Range("J32").Value = "Statement-non-opinion"
$ws.Range("I33").Value = "sd"
$ws.Range("J33").Value = "Statement-non-opinion"
$ws.Range("I51").Value = "sd"
$ws.Range("J51").Value = "Statement-non-opinion"
$ws.Range("I53").Value = "sd"
$ws.Range("J53").Value = "Statement-non-opinion"
$ws.Range("I60").Value = "qy"
$ws.Range("J60").Value = "Yes-No-Question"
$ws.Range("I69").Value = "sd"
$ws.Range("J69").Value = "Statement-non-opinion"
$ws.Range("I77").Value = "ba"
$ws.Range("J77").Value = "Appreciation"
$ws.Range("I81").Value = "b"
$ws.Range("J81").Value = "Acknowledge (Backchannel)"
$ws.Range("I89").Value = "sd"
$ws.Range("J89").Value = "Statement-non-opinion"
$ws.Range("I92").Value = "ba"
$ws.Range("J92").Value = "Appreciation"
$ws.Range("I94").Value = "sd"
$ws.Range("J94").Value = "Statement-non-opinion"
$ws.Range("I102").Value = "b"
$ws.Range("J102").Value = "Acknowledge (Backchannel)"
$ws.Range("I110").Value = "sv"
$ws.Range("J110").Value = "Statement-opinion"
$ws.Range("I122").Value = "aa"
$ws.Range("J122").Value = "Agree/Accept"
$ws.Range("I131").Value = "sv"
$ws.Range("J131").Value = "Statement-opinion"
$ws.Range("I139").Value = "sd"
$ws.Range("J139").Value = "Statement-non-opinion"
$ws.Range("I145").Value = "sv"
$ws.Range("J145").Value = "Statement-opinion"
$ws.Range("I149").Value = "aa"
$ws.Range("J149").Value = "Agree/Accept"
$ws.Range("I167").Value = "sd"
$ws.Range("J167").Value = "Statement-non-opinion"
$ws.Range("I181").Value = "ba"
$ws.Range("J181").Value = "Appreciation"
$ws.Range("I185").Value = "aa"
$ws.Range("J185").Value = "Agree/Accept"
$ws.Range("I188").Value = "sv"
$ws.Range("J188").Value = "Statement-opinion"
$ws.Range("I191").Value = "sd"
$ws.Range("J191").Value = "Statement-non-opinion"
$ws.Range("I206").Value = "ba"
$ws.Range("J206").Value = "Appreciation"
$ws.Range("I227").Value = "sv"
$ws.Range("J227").Value = "Statement-opinion"
$ws.Range("I240").Value = "%"
$ws.Range("J240").Value = "Uninterpretable"
$ws.Range("I247").Value = "sd"
$ws.Range("J247").Value = "Statement-non-opinion"
$ws.Range("I255").Value = "aa"
$ws.Range("J255").Value = "Agree/Accept"
$ws.Range("I265").Value = "sd"
$ws.Range("J265").Value = "Statement-non-opinion"
$ws.Range("I268").Value = "sd"
$ws.Range("J268").Value = "Statement-non-opinion"
$ws.Range("I269").Value = "aa"
$ws.Range("J269").Value = "Agree/Accept"
$ws.Range("I278").Value = "aa"
$ws.Range("J278").Value = "Agree/Accept"
$ws.Range("I288").Value = "sv"
$ws.Range("J288").Value = "Statement-opinion"
$ws.Range("I291").Value = "sd"
$ws.Range("J291").Value = "Statement-non-opinion"
$ws.Range("I301").Value = "aa"
$ws.Range("J301").Value = "Agree/Accept"
$ws.Range("I309").Value = "b"
$ws.Range("J309").Value = "Acknowledge (Backchannel)"
$ws.Range("I312").Value = "sv"
$ws.Range("J312").Value = "Statement-opinion"
$ws.Range("I327").Value = "sd"
$ws.Range("J327").Value = "Statement-non-opinion"
$ws.Range("I345").Value = "ba"
$ws.Range("J345").Value = "Appreciation"
$ws.Range("I361").Value = "b"
$ws.Range("J361").Value = "Acknowledge (Backchannel)"
$ws.Range("I372").Value = "b"
$ws.Range("J372").Value = "Acknowledge (Backchannel)"
$ws.Range("I374").Value = "sv"
$ws.Range("J374").Value = "Statement-opinion"
$ws.Range("I386").Value = "sd"
$ws.Range("J386").Value = "Statement-non-opinion"
$ws.Range("I387").Value = "sd"
$ws.Range("J387").Value = "Statement-non-opinion"
$ws.Range("I390").Value = "sd"
$ws.Range("J390").Value = "Statement-non-opinion"
$ws.Range("I396").Value = "aa"
$ws.Range("J396").Value = "Agree/Accept"
$ws.Range("I404").Value = "aa"
$ws.Range("J404").Value = "Agree/Accept"
$ws.Range("I410").Value = "%"
$ws.Range("J410").Value = "Uninterpretable"
$ws.Range("I419").Value = "sd"
$ws.Range("J419").Value = "Statement-non-opinion"
$ws.Range("I424").Value = "sd"
$ws.Range("J424").Value = "Statement-non-opinion"
$ws.Range("I427").Value = "sd"
$ws.Range("J427").Value = "Statement-non-opinion"
$ws.Range("I428").Value = "sd"
$ws.Range("J428").Value = "Statement-non-opinion"
$ws.Range("I433").Value = "sd"
$ws.Range("J433").Value = "Statement-non-opinion"
$ws.Range("I471").Value = "b"
$ws.Range("J471").Value = "Acknowledge (Backchannel)"
$ws.Range("I473").Value = "sd"
$ws.Range("J473").Value = "Statement-non-opinion"
$ws.Range("I486").Value = "sd"
$ws.Range("J486").Value = "Statement-non-opinion"
$ws.Range("I487").Value = "aa"
$ws.Range("J487").Value = "Agree/Accept"
$ws.Range("I492").Value = "%"
$ws.Range("J492").Value = "Uninterpretable"
$ws.Range("I500").Value = "sd"
$ws.Range("J500").Value = "Statement-non-opinion"
$ws.Range("I522").Value = "sd"
$ws.Range("J522").Value = "Statement-non-opinion"
$ws.Range("I535").Value = "sv"
$ws.Range("J535").Value = "Statement-opinion"
$ws.Range("I543").Value = "sv"
$ws.Range("J543").Value = "Statement-opinion"
$ws.Range("I561").Value = "sv"
$ws.Range("J561").Value = "Statement-opinion"
$ws.Range("I564").Value = "qy"
$ws.Range("J564").Value = "Yes-No-Question"
